$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = 1274.23779296875
$ws.Range("C2").Value = 0.9559
$ws.Range("D2").Value = 0.9169999957084656
$ws.Range("E2").Value = 1.233299970626831
$ws.Range("F2").Value = 0.8217999935150146
$ws.Range("H2").Value = 0.8628

# Row 3
$ws.Range("B3").Value = 1143.4365234375
$ws.Range("C3").Value = 0.9147
$ws.Range("D3").Value = 0.9035
$ws.Range("E3").Value = 1.264700055122375
$ws.Range("F3").Value = 0.8396999835968018
$ws.Range("H3").Value = 0.7429

# Row 4
$ws.Range("B4").Value = 758.06298828125
$ws.Range("C4").Value = 0.8856000000000001
$ws.Range("D4").Value = 0.887
$ws.Range("E4").Value = 0.9375
$ws.Range("F4").Value = 0.8205999732017517
$ws.Range("H4").Value = 0.5968

# Row 5
$ws.Range("B5").Value = 841.2050170898438
$ws.Range("C5").Value = 0.8873
$ws.Range("D5").Value = 0.8881
$ws.Range("E5").Value = 1.029199957847595
$ws.Range("F5").Value = 0.8338000178337097
$ws.Range("H5").Value = 0.6067

# Row 6
$ws.Range("B6").Value = 1155.6728515625
$ws.Range("C6").Value = 0.9143
$ws.Range("D6").Value = 0.9142
$ws.Range("E6").Value = 1.061100006103516
$ws.Range("F6").Value = 0.8205999732017517
$ws.Range("H6").Value = 0.8376

# Row 7 (D7 unchanged)
$ws.Range("B7").Value = 932.1348876953125
$ws.Range("C7").Value = 0.9387
$ws.Range("E7").Value = 1.124400019645691
$ws.Range("F7").Value = 0.876800000667572

# Row 8
$ws.Range("B8").Value = 1042.690795898438
$ws.Range("C8").Value = 0.9343
$ws.Range("D8").Value = 0.9282
$ws.Range("E8").Value = 1.157699942588806
$ws.Range("F8").Value = 0.883899986743927
$ws.Range("H8").Value = 0.9615

# Row 9
$ws.Range("B9").Value = 7147.44091796875
$ws.Range("C9").Value = 0.9211
$ws.Range("D9").Value = 0.9121
$ws.Range("E9").Value = 1.264700055122375
$ws.Range("F9").Value = 0.8205999732017517
$ws.Range("H9").Value = 5.5883
